# Extend the eval N/B cross table: add a "33" column (K) and a "33" row (12),
# mirroring the existing "32" column (J) / row (11) layout, plus fill in the
# newly playable eval5 data point (D12 = 4).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# --- New column K: header value + shaded (unused) cells below it, copied
#     from column J's formatting so the existing shared style is reused. ---
$ws.Range("K4").Value = 33

$ws.Range("J5:J11").Copy()
$ws.Range("K5:K11").PasteSpecial(-4122)
$ws.Range("K5:K11").ClearContents()

# --- New row 12: mirror row 11's layout for the shaded cells, then fill
#     in the real data (C12 = 33, D12 = 4). ---
$ws.Range("G11").Copy()
$ws.Range("G12").PasteSpecial(-4122)
$ws.Range("G12").ClearContents()

$ws.Range("J11").Copy()
$ws.Range("J12").PasteSpecial(-4122)
$ws.Range("J12").ClearContents()

$ws.Range("K11").Copy()
$ws.Range("K12").PasteSpecial(-4122)
$ws.Range("K12").ClearContents()

$ws.Range("C12").Value = 33
$ws.Range("D12").Value = 4

$excel.CutCopyMode = 0

$ws.Range("D13").Select()
